$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 74.10005433333333
$ws.Cells.Item(2, 8).Value = 222.300163
$ws.Cells.Item(2, 9).Value = 0.2282041889801584
$ws.Cells.Item(2, 10).Value = 0.2282041889801584
$ws.Cells.Item(2, 13).Value = 83.91225566666667
$ws.Cells.Item(2, 14).Value = 251.736767
$ws.Cells.Item(2, 15).Value = 0.9556261553553385
$ws.Cells.Item(2, 16).Value = 0.9556261553553385
$ws.Cells.Item(2, 17).Value = 6217.902704132558
$ws.Cells.Item(2, 18).Value = 55961.12433719302
$ws.Cells.Item(2, 19).Value = 0.2180778917510919
$ws.Cells.Item(2, 20).Value = 0.2180778917510919

$ws.Cells.Item(3, 7).Value = 74.10005433333333
$ws.Cells.Item(3, 8).Value = 222.300163
$ws.Cells.Item(3, 9).Value = 0.2282041889801584
$ws.Cells.Item(3, 10).Value = 0.2282041889801584
$ws.Cells.Item(3, 15).Value = 0.00439999103960854
$ws.Cells.Item(3, 16).Value = 0.00439999103960854
$ws.Cells.Item(3, 17).Value = 28.62909939208178
$ws.Cells.Item(3, 18).Value = 257.661894528736
$ws.Cells.Item(3, 19).Value = 0.001004096386713831
$ws.Cells.Item(3, 20).Value = 0.001004096386713831

$ws.Cells.Item(4, 7).Value = 74.10005433333333
$ws.Cells.Item(4, 8).Value = 222.300163
$ws.Cells.Item(4, 9).Value = 0.2282041889801584
$ws.Cells.Item(4, 10).Value = 0.2282041889801584
$ws.Cells.Item(4, 13).Value = 3.510050666666667
$ws.Cells.Item(4, 14).Value = 10.530152
$ws.Cells.Item(4, 15).Value = 0.03997385360505296
$ws.Cells.Item(4, 16).Value = 0.03997385360505297
$ws.Cells.Item(4, 17).Value = 260.0949451127529
$ws.Cells.Item(4, 18).Value = 2340.854506014776
$ws.Cells.Item(4, 19).Value = 0.009122200842352693
$ws.Cells.Item(4, 20).Value = 0.009122200842352693

$ws.Cells.Item(5, 9).Value = 0.566620969983319
$ws.Cells.Item(5, 10).Value = 0.566620969983319
$ws.Cells.Item(5, 13).Value = 83.91225566666667
$ws.Cells.Item(5, 14).Value = 251.736767
$ws.Cells.Item(5, 15).Value = 0.9556261553553385
$ws.Cells.Item(5, 16).Value = 0.9556261553553385
$ws.Cells.Item(5, 17).Value = 15438.77909175376
$ws.Cells.Item(5, 18).Value = 138949.0118257838
$ws.Cells.Item(5, 19).Value = 0.5414778190888718
$ws.Cells.Item(5, 20).Value = 0.5414778190888718

$ws.Cells.Item(6, 9).Value = 0.566620969983319
$ws.Cells.Item(6, 10).Value = 0.566620969983319
$ws.Cells.Item(6, 15).Value = 0.00439999103960854
$ws.Cells.Item(6, 16).Value = 0.00439999103960854
$ws.Cells.Item(6, 19).Value = 0.002493127190780903
$ws.Cells.Item(6, 20).Value = 0.002493127190780903

$ws.Cells.Item(7, 9).Value = 0.566620969983319
$ws.Cells.Item(7, 10).Value = 0.566620969983319
$ws.Cells.Item(7, 13).Value = 3.510050666666667
$ws.Cells.Item(7, 14).Value = 10.530152
$ws.Cells.Item(7, 15).Value = 0.03997385360505296
$ws.Cells.Item(7, 16).Value = 0.03997385360505297
$ws.Cells.Item(7, 17).Value = 645.8043156270019
$ws.Cells.Item(7, 18).Value = 5812.238840643017
$ws.Cells.Item(7, 19).Value = 0.0226500237036663
$ws.Cells.Item(7, 20).Value = 0.02265002370366631

$ws.Cells.Item(8, 7).Value = 66.62220766666667
$ws.Cells.Item(8, 8).Value = 199.866623
$ws.Cells.Item(8, 9).Value = 0.2051748410365226
$ws.Cells.Item(8, 10).Value = 0.2051748410365226
$ws.Cells.Item(8, 13).Value = 83.91225566666667
$ws.Cells.Item(8, 14).Value = 251.736767
$ws.Cells.Item(8, 15).Value = 0.9556261553553385
$ws.Cells.Item(8, 16).Value = 0.9556261553553385
$ws.Cells.Item(8, 17).Value = 5590.419722803093
$ws.Cells.Item(8, 18).Value = 50313.77750522784
$ws.Cells.Item(8, 19).Value = 0.1960704445153749
$ws.Cells.Item(8, 20).Value = 0.1960704445153748

$ws.Cells.Item(9, 7).Value = 66.62220766666667
$ws.Cells.Item(9, 8).Value = 199.866623
$ws.Cells.Item(9, 9).Value = 0.2051748410365226
$ws.Cells.Item(9, 10).Value = 0.2051748410365226
$ws.Cells.Item(9, 15).Value = 0.00439999103960854
$ws.Cells.Item(9, 16).Value = 0.00439999103960854
$ws.Cells.Item(9, 17).Value = 25.73997849487289
$ws.Cells.Item(9, 18).Value = 231.659806453856
$ws.Cells.Item(9, 19).Value = 0.0009027674621138063
$ws.Cells.Item(9, 20).Value = 0.0009027674621138062

$ws.Cells.Item(10, 7).Value = 66.62220766666667
$ws.Cells.Item(10, 8).Value = 199.866623
$ws.Cells.Item(10, 9).Value = 0.2051748410365226
$ws.Cells.Item(10, 10).Value = 0.2051748410365226
$ws.Cells.Item(10, 13).Value = 3.510050666666667
$ws.Cells.Item(10, 14).Value = 10.530152
$ws.Cells.Item(10, 15).Value = 0.03997385360505296
$ws.Cells.Item(10, 16).Value = 0.03997385360505297
$ws.Cells.Item(10, 17).Value = 233.8473244351885
$ws.Cells.Item(10, 18).Value = 2104.625919916696
$ws.Cells.Item(10, 19).Value = 0.008201629059033969
$ws.Cells.Item(10, 20).Value = 0.008201629059033969
